$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "30.082.89"
$ws.Cells.Item(2, 5).Value = "  -1.91%  "
$ws.Cells.Item(3, 4).Value = "2.104.65"
$ws.Cells.Item(3, 5).Value = "  -0.67%  "
$ws.Cells.Item(4, 5).Value = "  -0.62%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "345.71"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +2.08%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.61%  "
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.5183"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -1.73%  "
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.4440"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -2.41%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.09471"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +3.95%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "52.44"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -2.81%  "
$ws.Cells.Item(11, 5).Value = "  +0.28%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "25.28"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +3.31%  "
$ws.Cells.Item(13, 4).Value = "2.110.26"
$ws.Cells.Item(13, 5).Value = "  -0.51%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "6.736"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.49%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "8.131"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +0.45%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "99.79"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +1.23%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "0.00001168"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -0.19%  "
$ws.Cells.Item(18, 5).Value = "  -0.68%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "1.006"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.60%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "6.232"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -3.41%  "
$ws.Cells.Item(23, 4).Value = "30.167.50"
$ws.Cells.Item(23, 5).Value = "  -1.85%  "
$ws.Cells.Item(24, 5).Value = "  -1.93%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "2.331"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -2.02%  "
$ws.Cells.Item(26, 4).Value = "2.366.70"
$ws.Cells.Item(26, 5).Value = "  -0.14%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "22.07"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -1.93%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "164.50"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -0.72%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "2.551"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.01%  "
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "133.76"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -1.49%  "
$ws.Cells.Item(31, 5).Value = "  -3.14%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "0.1060"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -1.81%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "1.638"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.05%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "6.256"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -2.52%  "
$ws.Cells.Item(35, 5).Value = "  +0.30%  "
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "6.224"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "10.16"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -3.59%  "
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.02570"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -3.71%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.06791"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -1.30%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.2288"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -1.59%  "
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.6974"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +1.04%  "
$ws.Cells.Item(42, 5).Value = "  -0.57%  "
$ws.Cells.Item(43, 5).Value = "  +3.73%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.6711"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +3.51%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "14.29"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -5.67%  "
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "2.286"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -1.18%  "
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "3.642"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -1.69%  "
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "0.00000000358"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -2.65%  "
$ws.Cells.Item(49, 5).Value = "  -2.86%  "
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "82.65"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.50%  "
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.07207"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -1.47%  "
